$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the board names to add spaces around the hyphen separator
$ws.Range("A2").Value = "Andrew Shields - Java Foundations Project"
$ws.Range("A3").Value = "Moiya Josephs - Java Foundations Project"
$ws.Range("A4").Value = "Marielle Nolasco - .NET Foundations Project"

# Move the active selection to A4 (last edited cell)
$ws.Range("A4").Select()
